$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 413 (shifts existing rows 413+ down by one)
$ws.Rows("413:413").Insert()

# Populate the newly inserted row 413 with the new weekly price record
$ws.Range("A413").Value = 9
$ws.Range("B413").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C413").Value = "Metropolitana"
$ws.Range("D413").Value = 45194
$ws.Range("E413").Value = 13
$ws.Range("F413").Value = 100112043
$ws.Range("G413").Value = "Pepino ensalada"
$ws.Range("H413").Value = "Sin especificar"
$ws.Range("I413").Value = "Primera"
$ws.Range("J413").Value = 70
$ws.Range("K413").Value = 12000
$ws.Range("L413").Value = 13000
$ws.Range("M413").Value = 12500
$ws.Range("N413").Value = "$/caja 60 unidades"
$ws.Range("O413").Value = "Región de Arica y Parinacota"
$ws.Range("P413").Value = 208
$ws.Range("Q413").Value = 60
$ws.Range("R413").Value = "Hortaliza"
